$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.841.36"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.934.32"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "378.81"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.30%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "100.38"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("E11").Value = "  -0.31%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0850"
$c.ClearFormats()
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "3.394.82"
$ws.Range("E13").Value = "  -0.55%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "18.16"
$c.ClearFormats()
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("E15").Value = "  +71.31%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "7.60"
$c.ClearFormats()
$ws.Range("E16").Value = "  +3.74%  "
$ws.Range("D17").Value = "2.930.10"
$ws.Range("E17").Value = "  -6.29%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.992"
$c.ClearFormats()
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").Value = "50.793.28"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  -6.22%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.35"
$c.ClearFormats()
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("E22").Value = "  -0.23%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "69.40"
$c.ClearFormats()
$ws.Range("E23").Value = "  +1.70%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "266.24"
$c.ClearFormats()
$ws.Range("E24").Value = "  +2.15%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.ClearFormats()
$ws.Range("E25").Value = "  +13.83%  "
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("E27").Value = "  +0.00%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.05"
$c.ClearFormats()
$ws.Range("E28").Value = "  -7.40%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "25.51"
$c.ClearFormats()
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("E30").Value = "  -2.74%  "
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("E32").Value = "  +2.69%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "50.49"
$c.ClearFormats()
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +0.31%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "33.37"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  +4.78%  "
$ws.Range("E39").Value = "  +1.09%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "16.46"
$c.ClearFormats()
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.ClearFormats()
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "119.95"
$c.ClearFormats()
$ws.Range("E43").Value = "  -1.25%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "21.06"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("E45").Value = "  +7.08%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.01"
$c.ClearFormats()
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "2.006.59"
$ws.Range("E48").Value = "  +0.46%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.259"
$c.ClearFormats()
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("E50").Value = "  -4.75%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "5.27"
$c.ClearFormats()
$ws.Range("E51").Value = "  +4.68%  "
